# Applies targeted cell value corrections to the stock report sheet.
# Generated from the authoritative diff of F/G (qty/value) figures,
# two swapped row-pairs (B/E/F/G), and the resulting Sub Total / Grand Total rollups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 63
$ws.Range("G6").Value = 1882.44
$ws.Range("F7").Value = 94
$ws.Range("G7").Value = 4402.02
$ws.Range("B10").Value = 27096.88
$ws.Range("F70").Value = 5
$ws.Range("G70").Value = 674.75
$ws.Range("F77").Value = 243
$ws.Range("G77").Value = 11357.82
$ws.Range("F83").Value = 106
$ws.Range("G83").Value = 15971.02
$ws.Range("B90").Value = 171818.86
$ws.Range("F144").Value = 1002
$ws.Range("G144").Value = 8466.9
$ws.Range("F145").Value = 408
$ws.Range("G145").Value = 3259.92
$ws.Range("B147").Value = 13242.24
$ws.Range("F169").Value = 1
$ws.Range("G169").Value = 143.5
$ws.Range("B175").Value = 26852.35
$ws.Range("F270").Value = 14
$ws.Range("G270").Value = 451.36
$ws.Range("B275").Value = 5160.51
$ws.Range("F283").Value = 38
$ws.Range("G283").Value = 12975.86
$ws.Range("F288").Value = 37
$ws.Range("G288").Value = 3440.63
$ws.Range("F293").Value = 31
$ws.Range("G293").Value = 2179.92
$ws.Range("B304").Value = 169837.8
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("F455").Value = 42
$ws.Range("G455").Value = 2671.62
$ws.Range("B460").Value = 12597.96
$ws.Range("F485").Value = 9
$ws.Range("G485").Value = 1579.23
$ws.Range("B488").Value = 28985.92
$ws.Range("F599").Value = 1465
$ws.Range("G599").Value = 238956.15
$ws.Range("B606").Value = 394293.19
$ws.Range("F610").Value = 12
$ws.Range("G610").Value = 491.88
$ws.Range("F611").Value = 37
$ws.Range("G611").Value = 4802.23
$ws.Range("B618").Value = 42718.48
$ws.Range("B619").Value = 1655243.47
$ws.Range("B620").Value = 1655243.47
